# Apply updated "dSF" (column F) values per commit: "repull data, push all data, mean calculation"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8  = -3
    12 = -3
    13 = -4
    14 = -5
    16 = -3
    17 = -1
    19 = -1
    22 = 0
    23 = -3
    26 = -4
    31 = 2
    33 = -2
    36 = -3
    41 = -7
    42 = -2
    43 = -3
    44 = -5
    45 = 7
    54 = -2
    55 = 5
    56 = -4
    60 = 6
    63 = -7
    66 = 0
    67 = -1
    72 = -1
    74 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
